$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.476.42'
$ws.Range("E2").Value = '  +4.14%  '
$ws.Range("D3").Value = '2.989.44'
$ws.Range("E3").Value = '  +4.69%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.66'
$ws.Range("E5").Value = '  +8.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.63'
$ws.Range("E6").Value = '  +10.68%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.431'
$ws.Range("E8").Value = '  +7.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.53'
$ws.Range("E9").Value = '  +14.51%  '
$ws.Range("E10").Value = '  +12.51%  '
$ws.Range("E11").Value = '  +5.82%  '
$ws.Range("E12").Value = '  +5.73%  '
$ws.Range("D13").Value = '3.505.36'
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.41'
$ws.Range("E14").Value = '  +9.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000153'
$ws.Range("E15").Value = '  +15.48%  '
$ws.Range("D16").Value = '56.517.72'
$ws.Range("E16").Value = '  +4.34%  '
$ws.Range("D17").Value = '2.995.35'
$ws.Range("E17").Value = '  +5.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.86'
$ws.Range("E18").Value = '  +9.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.36'
$ws.Range("E19").Value = '  +9.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.78'
$ws.Range("E20").Value = '  +11.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.37'
$ws.Range("E21").Value = '  +10.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.477'
$ws.Range("E23").Value = '  +9.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.34'
$ws.Range("E24").Value = '  +7.16%  '
$ws.Range("E25").Value = '  +13.50%  '
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("D27").Value = '0.0₃0903'
$ws.Range("E27").Value = '  +14.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.54'
$ws.Range("E28").Value = '  +7.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.03'
$ws.Range("E29").Value = '  +14.54%  '
$ws.Range("E30").Value = '  +14.13%  '
$ws.Range("E31").Value = '  +11.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.54'
$ws.Range("E32").Value = '  +10.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.74'
$ws.Range("E33").Value = '  +13.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("E34").Value = '  +8.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +4.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.25'
$ws.Range("E36").Value = '  +4.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0674'
$ws.Range("E37").Value = '  +10.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.89'
$ws.Range("E38").Value = '  +4.67%  '
$ws.Range("D39").Value = '3.025.37'
$ws.Range("E39").Value = '  +5.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.89'
$ws.Range("E40").Value = '  +5.60%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.645'
$ws.Range("E42").Value = '  +7.80%  '
$ws.Range("D43").Value = '2.253.15'
$ws.Range("E43").Value = '  +11.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.40'
$ws.Range("E44").Value = '  +7.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.987'
$ws.Range("E45").Value = '  +6.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.59'
$ws.Range("E46").Value = '  +6.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.98'
$ws.Range("E47").Value = '  +26.11%  '
$ws.Range("E48").Value = '  +10.86%  '
$ws.Range("E49").Value = '  +8.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.98'
$ws.Range("E50").Value = '  +8.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0869'
$ws.Range("E51").Value = '  +10.61%  '
